# "Generate Report for Archive"
# The localization status report is regenerated: every "Ready for handoff"
# status cell moves to "In Translation" (the shared string is reused across
# the Overview roll-up sheet and the per-locale detail sheets), and the
# now-narrower status columns are resized to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E & F), rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 12.57

# --- zh-cn sheet: Status column (C), rows 2-4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsZhCn.Columns("C:C").ColumnWidth = 12.57

# --- de-de sheet: Status column (C), rows 2-4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"
$wsDeDe.Columns("C:C").ColumnWidth = 12.57
